# Add 2022-Q1 data.
#
# Current workbook layout:
#   Sheet1 "2021-Q4" (rId1) - fund detail table, already correct, untouched.
#   Sheet2 "总计"      (rId2) - 2-column summary table (date / count / value).
#
# Target workbook layout:
#   Sheet1 "2021-Q4" (rId1) - unchanged.
#   Sheet2 "2022-Q1" (rId2) - fund detail table for 2022-Q1 (reuses the
#                             "总计" sheet slot/r:id, content fully replaced).
#   Sheet3 "总计"      (rId3) - new sheet (copy of the old "总计" sheet, so it
#                             keeps the same header formatting), with a new
#                             row inserted on top for the 2022-Q1 summary.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(2)   # currently named "总计"

# 1. Duplicate the existing "总计" sheet *before* we overwrite it, so the
#    new summary sheet keeps identical header/row styling. The copy is
#    placed immediately after the source sheet.
$summarySheet.Copy($null, $summarySheet)
$newSummarySheet = $wb.Worksheets.Item(3)

# 2. Turn the duplicate into the updated "总计" sheet: push the existing
#    2021-Q4 summary row down to row 3 (copying A2's formatting onto A3
#    first, so the row-index cell keeps its bold/centered/bordered style),
#    and add the new 2022-Q1 row on top.
$newSummarySheet.Range("A2").Copy()
$newSummarySheet.Range("A3").PasteSpecial(-4122)

$newSummarySheet.Range("A3").Value = 1
$newSummarySheet.Range("B3").Value = "2021-Q4"
$newSummarySheet.Range("C3").Value = 4
$newSummarySheet.Range("D3").Value = 0.51

$newSummarySheet.Range("A2").Value = 0
$newSummarySheet.Range("B2").Value = "2022-Q1"
$newSummarySheet.Range("C2").Value = 4
$newSummarySheet.Range("D2").Value = 0.01

# 3. Repurpose the original "总计" sheet as the "2022-Q1" fund detail sheet.
# (rename it away from "总计" FIRST so the duplicate sheet can take that name)
$summarySheet.Name = "2022-Q1"
$newSummarySheet.Name = "总计"

# Extend the header style (bold/centered/bordered, same as B1:D1) across the
# new E1:H1 header cells.
$summarySheet.Range("D1").Copy()
$summarySheet.Range("E1:H1").PasteSpecial(-4122)

$summarySheet.Range("B1").Value = "基金代码"
$summarySheet.Range("C1").Value = "基金名称"
$summarySheet.Range("D1").Value = "基金规模"
$summarySheet.Range("E1").Value = "股票总仓位"
$summarySheet.Range("F1").Value = "仓位占比"
$summarySheet.Range("G1").Value = "持有市值(亿元)"
$summarySheet.Range("H1").Value = "仓位排名"

# Extend the row-index style (A2) down across A3:A5.
$summarySheet.Range("A2").Copy()
$summarySheet.Range("A3:A5").PasteSpecial(-4122)

# Fund code / scale / position columns are stored as plain text (to keep
# leading zeros like "009387" and trailing zeros like "4.30" intact), so
# format them as Text before writing the values.
$summarySheet.Range("B2:B5").NumberFormat = "@"
$summarySheet.Range("D2:G5").NumberFormat = "@"

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "009387"
$summarySheet.Range("C2").Value = "嘉实稳福混合A"
$summarySheet.Range("D2").Value = "0.08"
$summarySheet.Range("E2").Value = "34.71"
$summarySheet.Range("F2").Value = "4.30"
$summarySheet.Range("G2").Value = "0.0034"
$summarySheet.Range("H2").Value = 2

$summarySheet.Range("A3").Value = 1
$summarySheet.Range("B3").Value = "009649"
$summarySheet.Range("C3").Value = "嘉实精选平衡混合A"
$summarySheet.Range("D3").Value = "0.06"
$summarySheet.Range("E3").Value = "67.70"
$summarySheet.Range("F3").Value = "3.73"
$summarySheet.Range("G3").Value = "0.0022"
$summarySheet.Range("H3").Value = 9

$summarySheet.Range("A4").Value = 2
$summarySheet.Range("B4").Value = "009650"
$summarySheet.Range("C4").Value = "嘉实精选平衡混合C"
$summarySheet.Range("D4").Value = "0.01"
$summarySheet.Range("E4").Value = "67.70"
$summarySheet.Range("F4").Value = "3.73"
$summarySheet.Range("G4").Value = "0.0004"
$summarySheet.Range("H4").Value = 9

$summarySheet.Range("A5").Value = 3
$summarySheet.Range("B5").Value = "009388"
$summarySheet.Range("C5").Value = "嘉实稳福混合C"
$summarySheet.Range("D5").Value = "0.01"
$summarySheet.Range("E5").Value = "34.71"
$summarySheet.Range("F5").Value = "4.30"
$summarySheet.Range("G5").Value = "0.0004"
$summarySheet.Range("H5").Value = 2

# Restore the original active sheet (sheet copy/rename above shifts focus
# onto the newly created "总计" sheet).
$wb.Worksheets.Item(1).Activate()
